$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.061.12"
$ws.Range("D3").Value = "1.690.44"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.13"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.54%  "
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0886"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.43%  "
$ws.Range("D12").Value = "1.930.84"
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").Value = "1.692.33"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("E14").Value = "  -0.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.557"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "249.86"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.09%  "
$ws.Range("D18").Value = "28.032.97"
$ws.Range("E18").Value = "  +3.18%  "
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("E20").Value = "  -3.03%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("E28").Value = "  +0.42%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("E31").Value = "  +3.25%  "
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("E33").Value = "  -1.83%  "
$ws.Range("D34").Value = "1.450.89"
$ws.Range("E34").Value = "  -6.16%  "
$ws.Range("E35").Value = "  -2.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.948"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("E38").Value = "  -1.91%  "
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("E40").Value = "  -2.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.49"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.46%  "
$ws.Range("D44").Value = "1.837.10"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.24"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.798"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.73%  "
$ws.Range("E47").Value = "  +7.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "89.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("E50").Value = "  -0.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.89%  "
